$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.856.69'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').Value = '  +1.89%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.713.84'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').Value = '  +1.12%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '613.69'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').Value = '  +8.41%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '194.32'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').Value = '  +12.73%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.638'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').Value = '  +2.97%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').Value = '  +0.28%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.724'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').Value = '  +3.93%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').Value = '  +0.25%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '59.91'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').Value = '  +17.62%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000289'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').Value = '  -0.17%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.44'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').Value = '  +0.46%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.315.26'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').Value = '  +1.18%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.723.91'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').Value = '  +1.34%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.15'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').Value = '  +3.60%  '

$ws.Range('B17').Value = 'TRON'

$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.127'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').Value = '  +1.08%  '

$ws.Range('B18').Value = 'Chainlink'

$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.42'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').Value = '  +1.27%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.93'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').Value = '  +1.57%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.836.53'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').Value = '  +2.23%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '411.56'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').Value = '  +2.32%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.60'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').Value = '  +4.99%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '90.05'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').Value = '  +3.44%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.10'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').Value = '  +2.57%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.14'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').Value = '  +4.07%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.33'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').Value = '  +6.94%  '

$ws.Range('B27').Value = 'LEO'

$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.06'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').Value = '  +1.35%  '

$ws.Range('B28').Value = 'Toncoin'

$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.81'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').Value = '  +3.36%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.70'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').Value = '  +3.77%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.82'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').Value = '  +1.50%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.73'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').Value = '  +2.25%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.73'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').Value = '  +2.94%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '46.74'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').Value = '  +9.83%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.123'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').Value = '  +7.02%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '638.76'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').Value = '  +9.53%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '67.51'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').Value = '  +4.99%  '

$ws.Range('B37').Value = 'TheGraph'

$ws.Range('C37').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.415'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').Value = '  +5.74%  '

$ws.Range('B38').Value = 'PEPE'

$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0826'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').Value = '  -6.68%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').Value = '  -0.11%  '

$ws.Range('E40').Value = '  +0.08%  '

$ws.Range('E41').Value = '  +6.02%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.05'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').Value = '  +3.49%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0449'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').Value = '  +3.95%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.62'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').Value = '  +3.27%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.921.20'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').Value = '  +8.08%  '

$ws.Range('B46').Value = 'THORChain'

$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.36'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').Value = '  +2.78%  '

$ws.Range('B47').Value = 'Stellar'

$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.139'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').Value = '  +4.88%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.72'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').Value = '  +1.31%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '145.56'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').Value = '  +2.42%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.11'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').Value = '  -0.80%  '

$ws.Range('B51').Value = 'Stacks'

$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.78'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').Value = '  +2.43%  '
